$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 4194.273506371725
$ws.Range("C3").Value = 4066.502338727755
$ws.Range("C4").Value = 4066.502338727755
$ws.Range("C5").Value = 4066.502338727755
$ws.Range("C6").Value = 4066.502338727755
$ws.Range("C7").Value = 4003.881445601857
$ws.Range("C8").Value = 4003.881445601857
$ws.Range("C9").Value = 3885.184528486806
$ws.Range("C10").Value = 3885.184528486806
$ws.Range("C11").Value = 3885.184528486806
$ws.Range("C12").Value = 3885.184528486806
